$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new stockyard row (Y9) below the existing data, copying the
# formatting from the row directly above so it picks up the same styles
# (name column + numeric area column).
$ws.Range("A16:B16").Copy() | Out-Null
$ws.Range("A17:B17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("A17").Value = "Y9"
$ws.Range("B17").Value = 3200

# Relabel the header row from "Stock" / "area" to "name" / "area"
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "area"

# Append a second new stockyard row (Y80)
$ws.Range("A16:B16").Copy() | Out-Null
$ws.Range("A18:B18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("A18").Value = "Y80"
$ws.Range("B18").Value = 22681

# Match the saved selection/active-cell state
$ws.Range("F20").Select() | Out-Null

$wb.Save()
